$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 1 match name: MI vs CSK -> MI vs RCB
$ws.Range("C10").Value = "MI vs RCB"

# Fill in the rest of the contest match schedule (previously blank)
$ws.Range("C11").Value = "CSK vs DC"
$ws.Range("C12").Value = "SRH vs KKR"
$ws.Range("C13").Value = "RR vs PBKS"
$ws.Range("C14").Value = "KKR vs MI"
$ws.Range("C15").Value = "SRH vs RCB"
$ws.Range("C16").Value = "RR vs DC"
$ws.Range("C17").Value = "PBKS vs CSK"
$ws.Range("C18").Value = "MI vs SRH"
$ws.Range("C19").Value = "RCB vs KKR"
$ws.Range("C20").Value = "DC vs PBKS"
$ws.Range("C21").Value = "CSK vs RR"
$ws.Range("C22").Value = "DC vs MI"
$ws.Range("C23").Value = "PBKS vs SRH"
$ws.Range("C24").Value = "KKR vs CSK"

# Updated predicted scores for contest 1 (row 10)
$ws.Range("E10").Value = 20
$ws.Range("K10").Value = 0
$ws.Range("Q10").Value = 60
$ws.Range("T10").Value = 40
